# sp_Blitz Check ID List - v36, October 5 2014 public release update.
# Adds checks 131-151 (non-default database config options, Hekaton/In-Memory
# OLTP checks, database files on network shares/Azure, default-trace error and
# log-growth checks), adds the "Server Name" info row's URL, and refreshes the
# title row to the new version/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the title row (A1) to the new version/date.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "sp_Blitz® Check ID List - v36 Oct 5, 2014"

# ---------------------------------------------------------------------------
# 2. Row 195 ("Server Name" info row) gains a URL in column E.
# ---------------------------------------------------------------------------
$ws.Range("E195").Value = "http://BrentOzar.com/go/servername"

# ---------------------------------------------------------------------------
# 3. Append the new check rows (131-151) as rows 196-216.
#    Columns: A=CheckID, B=Priority, C=FindingsGroup, D=Finding, E=URL/info
# ---------------------------------------------------------------------------
$rows = @(
    ,@(196, 131, 210, 'Non-Default Database Config', 'Supplemental Logging Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(197, 132, 210, 'Non-Default Database Config', 'Snapshot Isolation Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(198, 133, 210, 'Non-Default Database Config', 'Read Committed Snapshot Isolation Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(199, 134, 210, 'Non-Default Database Config', 'Auto Create Stats Incremental Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(200, 135, 210, 'Non-Default Database Config', 'ANSI NULL Default Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(201, 136, 210, 'Non-Default Database Config', 'Recursive Triggers Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(202, 137, 210, 'Non-Default Database Config', 'Trustworthy Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(203, 138, 210, 'Non-Default Database Config', 'Forced Parameterization Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(204, 139, 210, 'Non-Default Database Config', 'Query Store Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(205, 140, 210, 'Non-Default Database Config', 'Change Data Capture Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(206, 141, 210, 'Non-Default Database Config', 'Containment Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(207, 142, 210, 'Non-Default Database Config', 'Target Recovery Time Changed', 'http://BrentOzar.com/go/dbdefaults')
    ,@(208, 143, 210, 'Non-Default Database Config', 'Delayed Durability Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(209, 144, 210, 'Non-Default Database Config', 'Memory Optimized Enabled', 'http://BrentOzar.com/go/dbdefaults')
    ,@(210, 145, 10, 'Performance', 'High Memory Use for In-Memory OLTP (Hekaton)', 'http://BrentOzar.com/go/hekaton')
    ,@(211, 146, 200, 'Performance', 'In-Memory OLTP (Hekaton) In Use', 'http://BrentOzar.com/go/hekaton')
    ,@(212, 147, 100, 'In-Memory OLTP (Hekaton)', 'Transaction Errors', 'http://BrentOzar.com/go/hekaton')
    ,@(213, 148, 50, 'Reliability', 'Database Files on Network File Shares', 'http://BrentOzar.com/go/nas')
    ,@(214, 149, 50, 'Reliability', 'Database Files Stored in Azure', 'http://BrentOzar.com/go/azurefiles')
    ,@(215, 150, 50, 'Reliability', 'Errors Logged Recently in the Default Trace', 'http://BrentOzar.com/go/defaulttrace')
    ,@(216, 151, 50, 'Performance', 'Log File Growths Slow', 'http://BrentOzar.com/go/filegrowth')
)

# Rows 196-214 have their E column as a plain (non-hyperlinked) URL; rows
# 210-214 get upgraded to real hyperlinks + the "Hyperlink" style below.
foreach ($r in $rows) {
    $rowNum = $r[0]
    $arr = New-Object 'object[,]' 1,5
    $arr[0,0] = $r[1]
    $arr[0,1] = $r[2]
    $arr[0,2] = $r[3]
    $arr[0,3] = $r[4]
    $arr[0,4] = $r[5]
    $ws.Range("A${rowNum}:E${rowNum}").Value = $arr
}

# ---------------------------------------------------------------------------
# 4. Checks 145-149 (rows 210-214) get real clickable hyperlinks in column E.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E210"), "http://BrentOzar.com/go/hekaton")
$ws.Range("E210").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E211"), "http://BrentOzar.com/go/hekaton")
$ws.Range("E211").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E212"), "http://BrentOzar.com/go/hekaton")
$ws.Range("E212").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E213"), "http://BrentOzar.com/go/nas")
$ws.Range("E213").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E214"), "http://BrentOzar.com/go/azurefiles")
$ws.Range("E214").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 5. Re-freeze the panes so the frozen header scrolls back up to the top
#    (topLeftCell B5) instead of staying scrolled down near the old last row.
# ---------------------------------------------------------------------------
$ws.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2").Select()

Write-Host "sp_Blitz Check ID List updated to v36 Oct 5, 2014 (checks 131-151 added)."
